$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n (column J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14: Average of SW(S*)/SW(OPT)  -> average of column N
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)  -> average of column Z
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)  -> min of column N
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)  -> max of column Z
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting for the new summary cells: bold 12pt font, vertically centered.
# Style B14 directly, then fan the resulting format out to B15:B17 via a
# format-only paste so only a single new cell style is interned.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the new summary rows
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Sheet view: scroll so row 10 is the top-left visible row, select A14:B17
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A14:B17").Select() | Out-Null

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
